$wb = $excel.ActiveWorkbook

# --- Sheet 1: IncomeReport ---
$ws1 = $wb.Worksheets.Item("IncomeReport")

$ws1.Range("A2").Value = "NET"
$ws1.Range("B2").Value = "212.17M -> 234.52M -> 253.86M -> 274.7M -> 290.18M"
$ws1.Range("C2").Value = 37
$ws1.Range("D2").Value = "-0.13 -> -0.2 -> -0.13 -> -0.14 -> -0.12"
$ws1.Range("E2").Value = 8
$ws1.Range("F2").Value = "(54.95M) -> 6.65M -> 6.13M -> 34.08M -> 20.81M"
$ws1.Range("G2").Value = 138
$ws1.Range("H2").Value = "5.75 <- N/A <- N/A <- N/A"
$ws1.Range("I2").Value = "20.63 <- 16.39 <- 21.85 <- 18.96"

$ws1.Columns.Item(2).ColumnWidth = 24.877604166666668
$ws1.Columns.Item(4).ColumnWidth = 19.877604166666668
$ws1.Columns.Item(6).ColumnWidth = 22.877604166666668
$ws1.Columns.Item(9).ColumnWidth = 15.877604166666666

# --- Sheet 2: EarningsReport ---
$ws2 = $wb.Worksheets.Item("EarningsReport")

$ws2.Range("A2").Value = "NET"
$ws2.Range("B2").Value = 0.03
$ws2.Range("C2").Value = 0.08
$ws2.Range("D2").Value = 166
$ws2.Range("E2").Value = "158, 157"
$ws2.Range("F2").Value = "30, 31"
$ws2.Range("G2").Value = 25
